# Refresh the "cryptos" price table (coinranking.com scrape).
# Each row is (rank, Coin, Link, Price, Volume(1h)); this pass updates the
# Price/Volume columns with the latest quote, and a few rows additionally
# swapped rank position (Coin/Link/Price/Volume all change together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.916.55'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '2.350.12'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = "'241.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = "'0.674"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.72%  '
$ws.Range("D7").Value = "'72.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.55%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("E10").Value = '  -2.86%  '
$ws.Range("D11").Value = "'58.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = "'33.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("D15").Value = '2.700.54'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("E16").Value = '  -5.07%  '
$ws.Range("D17").Value = "'0.906"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '2.354.67'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = '43.818.67'
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").Value = "'78.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = "'255.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("D24").Value = "'1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.45%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D27").Value = "'2.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.44%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = "'2.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.39%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = "'10.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("D30").Value = "'22.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("D31").Value = "'176.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = "'5.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.70%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = "'5.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("D38").Value = "'6.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("E39").Value = '  -4.66%  '
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").Value = "'67.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +23.88%  '
$ws.Range("D42").Value = "'5.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.52%  '
$ws.Range("D43").Value = "'0.109"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.11%  '
$ws.Range("D44").Value = "'9.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("E45").Value = '  +3.78%  '
$ws.Range("D46").Value = "'18.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.30%  '
$ws.Range("D47").Value = "'2.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").Value = "'99.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.17%  '
$ws.Range("E51").Value = '  -5.25%  '
